$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("miembros")

# Update all "fecha_inicio" (H column) values to 2024-01-01 (45292)
for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 8).Value = 45292
}

# Update condicion (F column) for a couple of members to "Inactivo"
$ws.Range("F6").Value = "Inactivo"
$ws.Range("F8").Value = "Inactivo"

# Update sheet view: remove the frozen/scrolled topLeftCell and change selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("I19").Select()
